$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 updates (odds movement) ---
$ws.Range("G4").Value = 1.45
$ws.Range("I4").Value = 7
$ws.Range("L4").Value = 7
$ws.Range("U4").Value = 2.1
$ws.Range("V4").Value = 1.67
$ws.Range("X4").Value = 6.5
$ws.Range("Y4").Value = 8.5
$ws.Range("Z4").Value = 9.5
$ws.Range("AB4").Value = 29
$ws.Range("AE4").Value = 21
$ws.Range("AO4").Value = 7.5
$ws.Range("AQ4").Value = 23
$ws.Range("AU4").Value = 9.5
$ws.Range("AW4").Value = 8
$ws.Range("BB4").Value = 401

# --- Row 5 updates (odds movement) ---
$ws.Range("O5").Value = 1.33
$ws.Range("P5").Value = 3.4
$ws.Range("Q5").Value = 2.06
$ws.Range("R5").Value = 1.84

# --- Insert two new match rows at position 9 (pushes old rows 9-11 down to 11-13) ---
$ws.Range("A9:A10").EntireRow.Insert()

# --- Row 9: New match: Cuiaba vs Flamengo RJ ---
$ws.Range("A9").Value = "Q9lyYnY0"
$ws.Range("B9").Value = "20/11/2024"
$ws.Range("C9").Value = "19:00"
$ws.Range("D9").Value = "BRAZIL - SERIE A BETANO"
$ws.Range("E9").Value = "Cuiaba"
$ws.Range("F9").Value = "Flamengo RJ"
$ws.Range("G9").Value = 3.75
$ws.Range("H9").Value = 3.4
$ws.Range("I9").Value = 2
$ws.Range("J9").Value = 4.5
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 2.75
$ws.Range("M9").Value = 1.08
$ws.Range("N9").Value = 8
$ws.Range("O9").Value = 1.4
$ws.Range("P9").Value = 3
$ws.Range("Q9").Value = 2.25
$ws.Range("R9").Value = 1.62
$ws.Range("S9").Value = 1.5
$ws.Range("T9").Value = 2.5
$ws.Range("U9").Value = 2
$ws.Range("V9").Value = 1.75
$ws.Range("W9").Value = 9
$ws.Range("X9").Value = 17
$ws.Range("Y9").Value = 13
$ws.Range("Z9").Value = 41
$ws.Range("AA9").Value = 34
$ws.Range("AB9").Value = 41
$ws.Range("AC9").Value = 8
$ws.Range("AD9").Value = 6.5
$ws.Range("AE9").Value = 19
$ws.Range("AF9").Value = 67
$ws.Range("AG9").Value = 6
$ws.Range("AH9").Value = 8.5
$ws.Range("AI9").Value = 9
$ws.Range("AJ9").Value = 17
$ws.Range("AK9").Value = 19
$ws.Range("AL9").Value = 34
$ws.Range("AM9").Value = 451
$ws.Range("AN9").Value = 5.5
$ws.Range("AO9").Value = 23
$ws.Range("AP9").Value = 34
$ws.Range("AQ9").Value = 81
$ws.Range("AR9").Value = 126
$ws.Range("AS9").Value = 301
$ws.Range("AT9").Value = 2.5
$ws.Range("AU9").Value = 8.5
$ws.Range("AV9").Value = 67
$ws.Range("AW9").Value = 4
$ws.Range("AX9").Value = 11
$ws.Range("AY9").Value = 23
$ws.Range("AZ9").Value = 41
$ws.Range("BA9").Value = 67
$ws.Range("BB9").Value = 201
$ws.Range("BC9").Value = 126
$ws.Range("BD9").Value = 126

# --- Row 10: New match: Gremio vs Juventude ---
$ws.Range("A10").Value = "n319JmJ7"
$ws.Range("B10").Value = "20/11/2024"
$ws.Range("C10").Value = "19:00"
$ws.Range("D10").Value = "BRAZIL - SERIE A BETANO"
$ws.Range("E10").Value = "Gremio"
$ws.Range("F10").Value = "Juventude"
$ws.Range("G10").Value = 1.73
$ws.Range("H10").Value = 3.6
$ws.Range("I10").Value = 5.25
$ws.Range("J10").Value = 2.3
$ws.Range("K10").Value = 2.2
$ws.Range("L10").Value = 5
$ws.Range("M10").Value = 1.06
$ws.Range("N10").Value = 10
$ws.Range("O10").Value = 1.29
$ws.Range("P10").Value = 3.75
$ws.Range("Q10").Value = 1.95
$ws.Range("R10").Value = 1.95
$ws.Range("S10").Value = 1.4
$ws.Range("T10").Value = 2.75
$ws.Range("U10").Value = 1.8
$ws.Range("V10").Value = 1.95
$ws.Range("W10").Value = 7
$ws.Range("X10").Value = 8
$ws.Range("Y10").Value = 8.5
$ws.Range("Z10").Value = 13
$ws.Range("AA10").Value = 13
$ws.Range("AB10").Value = 26
$ws.Range("AC10").Value = 10
$ws.Range("AD10").Value = 7
$ws.Range("AE10").Value = 15
$ws.Range("AF10").Value = 51
$ws.Range("AG10").Value = 15
$ws.Range("AH10").Value = 26
$ws.Range("AI10").Value = 17
$ws.Range("AJ10").Value = 51
$ws.Range("AK10").Value = 41
$ws.Range("AL10").Value = 41
$ws.Range("AM10").Value = 251
$ws.Range("AN10").Value = 3.6
$ws.Range("AO10").Value = 9
$ws.Range("AP10").Value = 21
$ws.Range("AQ10").Value = 29
$ws.Range("AR10").Value = 51
$ws.Range("AS10").Value = 151
$ws.Range("AT10").Value = 2.75
$ws.Range("AU10").Value = 8.5
$ws.Range("AV10").Value = 51
$ws.Range("AW10").Value = 6.5
$ws.Range("AX10").Value = 26
$ws.Range("AY10").Value = 34
$ws.Range("AZ10").Value = 81
$ws.Range("BA10").Value = 101
$ws.Range("BB10").Value = 251
$ws.Range("BC10").Value = 126
$ws.Range("BD10").Value = 126

# --- Row 11: LATVIA - VIRSLIGA: Grobina vs Alberts JDFS (odds updated) ---
$ws.Range("A11").Value = "468pA9I6"
$ws.Range("B11").Value = "20/11/2024"
$ws.Range("C11").Value = "08:00"
$ws.Range("D11").Value = "LATVIA - VIRSLIGA"
$ws.Range("E11").Value = "Grobina"
$ws.Range("F11").Value = "Alberts JDFS"
$ws.Range("G11").Value = 1.35
$ws.Range("H11").Value = 4.25
$ws.Range("I11").Value = 7.3
$ws.Range("J11").Value = 1.85
$ws.Range("K11").Value = 2.3
$ws.Range("L11").Value = 6.6
$ws.Range("M11").Value = 1.04
$ws.Range("N11").Value = 10
$ws.Range("O11").Value = 1.19
$ws.Range("P11").Value = 3.65
$ws.Range("Q11").Value = 1.65
$ws.Range("R11").Value = 1.98
$ws.Range("S11").Value = 1.32
$ws.Range("T11").Value = 3.1
$ws.Range("U11").Value = 2
$ws.Range("V11").Value = 1.77
$ws.Range("W11").Value = 5.8
$ws.Range("X11").Value = 5.4
$ws.Range("Y11").Value = 7.1
$ws.Range("Z11").Value = 7.2
$ws.Range("AA11").Value = 9.5
$ws.Range("AB11").Value = 22
$ws.Range("AC11").Value = 11.5
$ws.Range("AD11").Value = 7.5
$ws.Range("AE11").Value = 16.5
$ws.Range("AF11").Value = 70
$ws.Range("AG11").Value = 16
$ws.Range("AH11").Value = 40
$ws.Range("AI11").Value = 18.5
$ws.Range("AJ11").Value = 120
$ws.Range("AK11").Value = 65
$ws.Range("AL11").Value = 55
$ws.Range("AM11").Value = 500
$ws.Range("AN11").Value = 3.1
$ws.Range("AO11").Value = 6.2
$ws.Range("AP11").Value = 17
$ws.Range("AQ11").Value = 17.5
$ws.Range("AR11").Value = 50
$ws.Range("AS11").Value = 250
$ws.Range("AT11").Value = 2.87
$ws.Range("AU11").Value = 8.5
$ws.Range("AV11").Value = 80
$ws.Range("AW11").Value = 8.5
$ws.Range("AX11").Value = 45
$ws.Range("AY11").Value = 45
$ws.Range("AZ11").Value = 300
$ws.Range("BA11").Value = 350
$ws.Range("BB11").Value = 400
$ws.Range("BC11").Value = 51
$ws.Range("BD11").Value = 51

# --- Row 13: WALES - CYMRU PREMIER: Briton Ferry vs TNS (odds updated) ---
$ws.Range("A13").Value = "6uOnIaCm"
$ws.Range("B13").Value = "20/11/2024"
$ws.Range("C13").Value = "16:45"
$ws.Range("D13").Value = "WALES - CYMRU PREMIER"
$ws.Range("E13").Value = "Briton Ferry"
$ws.Range("F13").Value = "TNS"
$ws.Range("G13").Value = 35
$ws.Range("H13").Value = 8.25
$ws.Range("I13").Value = 1.05
$ws.Range("J13").Value = 23
$ws.Range("K13").Value = 3.35
$ws.Range("L13").Value = 1.3
$ws.Range("M13").Value = 1.02
$ws.Range("N13").Value = 10
$ws.Range("O13").Value = 1.08
$ws.Range("P13").Value = 6.4
$ws.Range("Q13").Value = 1.28
$ws.Range("R13").Value = 3.35
$ws.Range("S13").Value = 1.2
$ws.Range("T13").Value = 4.1
$ws.Range("U13").Value = 2.77
$ws.Range("V13").Value = 1.39
$ws.Range("W13").Value = 150
$ws.Range("X13").Value = 800
$ws.Range("Y13").Value = 200
$ws.Range("Z13").Value = 1000
$ws.Range("AA13").Value = 500
$ws.Range("AB13").Value = 700
$ws.Range("AC13").Value = 21
$ws.Range("AD13").Value = 27
$ws.Range("AE13").Value = 70
$ws.Range("AF13").Value = 400
$ws.Range("AG13").Value = 10.75
$ws.Range("AH13").Value = 6.5
$ws.Range("AI13").Value = 15.5
$ws.Range("AJ13").Value = 5.7
$ws.Range("AK13").Value = 13
$ws.Range("AL13").Value = 60
$ws.Range("AM13").Value = 900
$ws.Range("AN13").Value = 30
$ws.Range("AO13").Value = 350
$ws.Range("AP13").Value = 175
$ws.Range("AR13").Value = 500
$ws.Range("AT13").Value = 4.1
$ws.Range("AU13").Value = 14.5
$ws.Range("AV13").Value = 150
$ws.Range("AW13").Value = 3
$ws.Range("AX13").Value = 3.85
$ws.Range("AY13").Value = 18
$ws.Range("AZ13").Value = 7.1
$ws.Range("BA13").Value = 32
$ws.Range("BB13").Value = 300

